$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the generated-on timestamp in A1
$ws.Range("A1").Value = "Reporte generado el 18/05/2025 a las 21:00"

# Update row 10 (item 8)
$ws.Range("B10").Value = "Hola"
$ws.Range("C10").Value = 21
$ws.Range("D10").Value = 123

# Update row 11 (item 9)
$ws.Range("B11").Value = "Hola 123"
$ws.Range("D11").Value = 12
$ws.Range("E11").Value = "unidad"

# Remove rows 12-17 (items 10-15), which are no longer part of the report
$ws.Rows("12:17").Delete()
